$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1347.7273
$ws.Range("J17").Value = 1353.8966
$ws.Range("L17").Value = 4061.6898
$ws.Range("N17").Value = -4397.6898
$ws.Range("H33").Value = 571.2353000000001
$ws.Range("I33").Value = 347
$ws.Range("J33").Value = 1300
$ws.Range("K33").Value = 347
$ws.Range("L33").Value = 1300
$ws.Range("M33").Value = -118
$ws.Range("N33").Value = -1758
$ws.Range("H40").Value = 5991.5
$ws.Range("I40").Value = 5990
$ws.Range("K40").Value = 5990
$ws.Range("M40").Value = -5815
$ws.Range("H41").Value = 591.5
$ws.Range("J41").Value = 786.36365
$ws.Range("L41").Value = 786.36365
$ws.Range("N41").Value = -1666.36365
$ws.Range("H64").Value = 53666.332
$ws.Range("I64").Value = 101582.664
$ws.Range("J64").Value = 5750
$ws.Range("K64").Value = 101582.664
$ws.Range("L64").Value = 5750
$ws.Range("M64").Value = -101334.664
$ws.Range("N64").Value = -6246
$ws.Range("H67").Value = 53666.332
$ws.Range("I67").Value = 101582.664
$ws.Range("J67").Value = 5750
$ws.Range("K67").Value = 101582.664
$ws.Range("L67").Value = 5750
$ws.Range("M67").Value = -100724.664
$ws.Range("N67").Value = -7466
$ws.Range("H80").Value = 251339.38
$ws.Range("I80").Value = 667157
$ws.Range("J80").Value = 1848.8
$ws.Range("K80").Value = 2001471
$ws.Range("L80").Value = 5546.4
$ws.Range("M80").Value = -2000473
$ws.Range("N80").Value = -7542.4
$ws.Range("H83").Value = 251339.38
$ws.Range("I83").Value = 667157
$ws.Range("J83").Value = 1848.8
$ws.Range("K83").Value = 6004413
$ws.Range("L83").Value = 16639.2
$ws.Range("M83").Value = -5999421
$ws.Range("N83").Value = -26623.2
$ws.Range("H97").Value = 2443.25
$ws.Range("J97").Value = 2078
$ws.Range("L97").Value = 6234
$ws.Range("N97").Value = -7226
$ws.Range("H103").Value = 874014.3
$ws.Range("I103").Value = 2038361.5
$ws.Range("J103").Value = 754
$ws.Range("K103").Value = 6115084.5
$ws.Range("L103").Value = 2262
$ws.Range("M103").Value = -6114498.5
$ws.Range("N103").Value = -3434
$ws.Range("H134").Value = 92719.53
$ws.Range("J134").Value = 92719.53
$ws.Range("L134").Value = 92719.53
$ws.Range("N134").Value = -102859.53
$ws.Range("H137").Value = 569130.6
$ws.Range("I137").Value = 774806.3
$ws.Range("K137").Value = 2324418.9
$ws.Range("M137").Value = -2321868.9
$ws.Range("H138").Value = 138070.75
$ws.Range("J138").Value = 4299.88
$ws.Range("L138").Value = 12899.64
$ws.Range("N138").Value = -23179.64
$ws.Range("H141").Value = 8158.05
$ws.Range("I141").Value = 7945.067
$ws.Range("K141").Value = 23835.201
$ws.Range("M141").Value = -18655.201

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3536.5
$ws.Range("I2").Value = 3910.6316
$ws.Range("K2").Value = 3910.6316
$ws.Range("M2").Value = -3797.6316
$ws.Range("H32").Value = 617782.7
$ws.Range("I32").Value = 642260.4
$ws.Range("J32").Value = 14000
$ws.Range("K32").Value = 642260.4
$ws.Range("L32").Value = 14000
$ws.Range("M32").Value = -641973.4
$ws.Range("N32").Value = -14574
$ws.Range("H45").Value = 73475.24000000001
$ws.Range("I45").Value = 90889.086
$ws.Range("K45").Value = 90889.086
$ws.Range("M45").Value = -90512.086
$ws.Range("H116").Value = 3536.5
$ws.Range("I116").Value = 3910.6316
$ws.Range("K116").Value = 3910.6316
$ws.Range("M116").Value = -1616.6316
$ws.Range("H122").Value = 1435079.6
$ws.Range("I122").Value = 6649.8
$ws.Range("K122").Value = 19949.4
$ws.Range("M122").Value = -17499.4
$ws.Range("H132").Value = 2248.8206
$ws.Range("I132").Value = 1396.7931
$ws.Range("J132").Value = 4719.7
$ws.Range("K132").Value = 4190.379300000001
$ws.Range("L132").Value = 14159.1
$ws.Range("M132").Value = -1660.379300000001
$ws.Range("N132").Value = -19219.1

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3536.5
$ws.Range("I3").Value = 3910.6316
$ws.Range("K3").Value = 3910.6316
$ws.Range("M3").Value = -3796.6316
$ws.Range("H99").Value = 37445.453
$ws.Range("I99").Value = 60835.332
$ws.Range("K99").Value = 60835.332
$ws.Range("M99").Value = -59337.332

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2967.3
$ws.Range("I31").Value = 2116.36
$ws.Range("K31").Value = 2116.36
$ws.Range("M31").Value = -1821.36
$ws.Range("H34").Value = 2967.3
$ws.Range("I34").Value = 2116.36
$ws.Range("K34").Value = 2116.36
$ws.Range("M34").Value = -1914.36
$ws.Range("H58").Value = 3194.9443
$ws.Range("I58").Value = 3078.0908
$ws.Range("J58").Value = 3378.5715
$ws.Range("K58").Value = 3078.0908
$ws.Range("L58").Value = 3378.5715
$ws.Range("M58").Value = -2875.0908
$ws.Range("N58").Value = -3784.5715
$ws.Range("H105").Value = 113492.52
$ws.Range("I105").Value = 126463.12
$ws.Range("K105").Value = 126463.12
$ws.Range("M105").Value = -124716.12
$ws.Range("H134").Value = 2072.6553
$ws.Range("I134").Value = 2118.1072
$ws.Range("K134").Value = 6354.321599999999
$ws.Range("M134").Value = -3819.321599999999
$ws.Range("H136").Value = 3194.9443
$ws.Range("I136").Value = 3078.0908
$ws.Range("J136").Value = 3378.5715
$ws.Range("K136").Value = 9234.2724
$ws.Range("L136").Value = 10135.7145
$ws.Range("M136").Value = -6684.2724
$ws.Range("N136").Value = -15235.7145

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 18334172
$ws.Range("J23").Value = 20371202
$ws.Range("L23").Value = 61113606
$ws.Range("N23").Value = -61114076
$ws.Range("H55").Value = 11817
$ws.Range("J55").Value = 12155.735
$ws.Range("L55").Value = 36467.205
$ws.Range("N55").Value = -36821.205
$ws.Range("H97").Value = 77624.375
$ws.Range("I97").Value = 100115.836
$ws.Range("J97").Value = 10150
$ws.Range("K97").Value = 300347.508
$ws.Range("L97").Value = 30450
$ws.Range("M97").Value = -299851.508
$ws.Range("N97").Value = -31442
$ws.Range("H137").Value = 8855.6
$ws.Range("J137").Value = 10419.5
$ws.Range("L137").Value = 31258.5
$ws.Range("N137").Value = -41458.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 15769.692
$ws.Range("I80").Value = 20502.75
$ws.Range("J80").Value = 8196.799999999999
$ws.Range("K80").Value = 20502.75
$ws.Range("L80").Value = 8196.799999999999
$ws.Range("M80").Value = -19504.75
$ws.Range("N80").Value = -10192.8
$ws.Range("H83").Value = 15769.692
$ws.Range("I83").Value = 20502.75
$ws.Range("J83").Value = 8196.799999999999
$ws.Range("K83").Value = 102513.75
$ws.Range("L83").Value = 40984
$ws.Range("M83").Value = -97521.75
$ws.Range("N83").Value = -50968
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H102").Value = 8990
$ws.Range("I102").Value = 9655
$ws.Range("K102").Value = 9655
$ws.Range("M102").Value = -8033

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 17131.223
$ws.Range("I7").Value = 20858.73
$ws.Range("J7").Value = 7439.7
$ws.Range("K7").Value = 20858.73
$ws.Range("L7").Value = 7439.7
$ws.Range("M7").Value = -20746.73
$ws.Range("N7").Value = -7663.7
$ws.Range("H16").Value = 5827.75
$ws.Range("I16").Value = 5793.1
$ws.Range("J16").Value = 6001
$ws.Range("K16").Value = 5793.1
$ws.Range("L16").Value = 6001
$ws.Range("M16").Value = -5623.1
$ws.Range("N16").Value = -6341
$ws.Range("H40").Value = 34121.684
$ws.Range("I40").Value = 55553.625
$ws.Range("K40").Value = 55553.625
$ws.Range("M40").Value = -55417.625
$ws.Range("H122").Value = 10777.444
$ws.Range("I122").Value = 22498.5
$ws.Range("J122").Value = 7428.5713
$ws.Range("K122").Value = 67495.5
$ws.Range("L122").Value = 22285.7139
$ws.Range("M122").Value = -65045.5
$ws.Range("N122").Value = -27185.7139
$ws.Range("H126").Value = 17131.223
$ws.Range("I126").Value = 20858.73
$ws.Range("J126").Value = 7439.7
$ws.Range("K126").Value = 62576.19
$ws.Range("L126").Value = 22319.1
$ws.Range("M126").Value = -60106.19
$ws.Range("N126").Value = -27259.1
$ws.Range("H132").Value = 556447.0600000001
$ws.Range("J132").Value = 4404.636
$ws.Range("L132").Value = 13213.908
$ws.Range("N132").Value = -18273.908

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 7710.7144
$ws.Range("I122").Value = 6816.6665
$ws.Range("J122").Value = 7954.5454
$ws.Range("K122").Value = 20449.9995
$ws.Range("L122").Value = 23863.6362
$ws.Range("M122").Value = -17999.9995
$ws.Range("N122").Value = -28763.6362
$ws.Range("H132").Value = 3586.6077
$ws.Range("I132").Value = 3496.0657
$ws.Range("J132").Value = 3893.4443
$ws.Range("K132").Value = 10488.1971
$ws.Range("L132").Value = 11680.3329
$ws.Range("M132").Value = -7958.197100000001
$ws.Range("N132").Value = -16740.3329
